$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column L width (col 12): 25 -> 26
$ws.Columns.Item(12).ColumnWidth = 26

# Row 2
$ws.Range("E2").Value = '2026-02-06 00:48:02'
$ws.Range("G2").Value = '108 cm'
$ws.Range("H2").Value = '99%'
$ws.Range("I2").Value = '0.0 mm'
$ws.Range("K2").Value = '0.0 MJ/m2'
$ws.Range("M2").Value = '-0.6 °C 0:12 TU'
$ws.Range("N2").Value = '-0.7 °C 0:28 TU'
$ws.Range("O2").Value = '-0.6 °C'

# Row 3
$ws.Range("E3").Value = '2026-02-06 00:48:04'
$ws.Range("G3").Value = '175 cm'
$ws.Range("H3").Value = '92%'
$ws.Range("I3").Value = '0.0 mm'
$ws.Range("K3").Value = '0.0 MJ/m2'
$ws.Range("L3").Value = '32.0 km/h - 228º 0:22 TU'
$ws.Range("M3").Value = '-2.1 °C 0:23 TU'
$ws.Range("N3").Value = '-3.5 °C 0:11 TU'
$ws.Range("O3").Value = '-2.8 °C'

# Row 4
$ws.Range("E4").Value = '2026-02-06 00:48:07'
$ws.Range("H4").Value = '50%'
$ws.Range("I4").Value = '0.0 mm'
$ws.Range("J4").Value = '990.6 hPa'
$ws.Range("K4").Value = '0.0 MJ/m2'
$ws.Range("L4").Value = '37.8 km/h - 291º 0:21 TU'
$ws.Range("M4").Value = '15.0 °C 0:02 TU'
$ws.Range("N4").Value = '14.5 °C 0:28 TU'
$ws.Range("O4").Value = '14.8 °C'

# Row 5
$ws.Range("E5").Value = '2026-02-06 00:48:10'
$ws.Range("H5").Value = '69%'
$ws.Range("I5").Value = '0.0 mm'
$ws.Range("J5").Value = '991.3 hPa'
$ws.Range("K5").Value = '0.0 MJ/m2'
$ws.Range("L5").Value = '9.7 km/h - 326º 0:00 TU'
$ws.Range("M5").Value = '10.2 °C 0:21 TU'
$ws.Range("N5").Value = '9.2 °C 0:06 TU'
$ws.Range("O5").Value = '9.8 °C'

# Row 6
$ws.Range("E6").Value = '2026-02-06 00:48:13'
$ws.Range("H6").Value = '51%'
$ws.Range("I6").Value = '0.0 mm'
$ws.Range("J6").Value = '993.0 hPa'
$ws.Range("K6").Value = '0.0 MJ/m2'
$ws.Range("L6").Value = '25.2 km/h - 284º 0:24 TU'
$ws.Range("M6").Value = '15.3 °C 0:01 TU'
$ws.Range("N6").Value = '14.5 °C 0:23 TU'
$ws.Range("O6").Value = '14.8 °C'

# Row 7
$ws.Range("E7").Value = '2026-02-06 00:48:15'
$ws.Range("H7").Value = '67%'
$ws.Range("I7").Value = '0.0 mm'
$ws.Range("J7").Value = '992.8 hPa'
$ws.Range("K7").Value = '0.0 MJ/m2'
$ws.Range("L7").Value = '31.7 km/h - 241º 0:22 TU'
$ws.Range("M7").Value = '10.5 °C 0:08 TU'
$ws.Range("N7").Value = '10.2 °C 0:23 TU'
$ws.Range("O7").Value = '10.4 °C'

# Row 8
$ws.Range("E8").Value = '2026-02-06 00:48:18'
$ws.Range("H8").Value = '86%'
$ws.Range("I8").Value = '0.0 mm'
$ws.Range("K8").Value = '0.0 MJ/m2'
$ws.Range("L8").Value = '7.2 km/h - 92º 0:06 TU'
$ws.Range("M8").Value = '8.6 °C 0:21 TU'
$ws.Range("N8").Value = '7.4 °C 0:08 TU'
$ws.Range("O8").Value = '7.9 °C'

# Row 9
$ws.Range("E9").Value = '2026-02-06 00:48:20'
$ws.Range("H9").Value = '100%'
$ws.Range("I9").Value = '0.0 mm'
$ws.Range("M9").Value = '2.9 °C 0:24 TU'
$ws.Range("N9").Value = '2.2 °C 0:18 TU'
$ws.Range("O9").Value = '2.5 °C'

# Row 10
$ws.Range("E10").Value = '2026-02-06 00:48:23'
$ws.Range("H10").Value = '100%'
$ws.Range("I10").Value = '0.0 mm'
$ws.Range("M10").Value = '6.0 °C 0:24 TU'
$ws.Range("N10").Value = '5.0 °C 0:13 TU'
$ws.Range("O10").Value = '5.5 °C'

# Row 11
$ws.Range("E11").Value = '2026-02-06 00:48:26'
$ws.Range("G11").Value = '1 cm'
$ws.Range("H11").Value = '84%'
$ws.Range("I11").Value = '0.0 mm'
$ws.Range("J11").Value = '994.2 hPa'
$ws.Range("K11").Value = '0.0 MJ/m2'
$ws.Range("L11").Value = '28.8 km/h - 187º 0:06 TU'
$ws.Range("M11").Value = '5.2 °C 0:00 TU'
$ws.Range("N11").Value = '4.8 °C 0:29 TU'
$ws.Range("O11").Value = '5.0 °C'

# Row 12
$ws.Range("E12").Value = '2026-02-06 00:48:28'
$ws.Range("H12").Value = '58%'
$ws.Range("I12").Value = '0.0 mm'
$ws.Range("K12").Value = '0.0 MJ/m2'
$ws.Range("L12").Value = '24.5 km/h - 226º 0:27 TU'
$ws.Range("M12").Value = '15.2 °C 0:02 TU'
$ws.Range("N12").Value = '12.6 °C 0:26 TU'
$ws.Range("O12").Value = '13.7 °C'

# Row 13
$ws.Range("E13").Value = '2026-02-06 00:48:31'
$ws.Range("H13").Value = '86%'
$ws.Range("I13").Value = '0.0 mm'
$ws.Range("M13").Value = '8.3 °C 0:29 TU'
$ws.Range("N13").Value = '6.2 °C 0:00 TU'
$ws.Range("O13").Value = '7.5 °C'

# Row 14
$ws.Range("E14").Value = '2026-02-06 00:48:33'
$ws.Range("G14").Value = '76 cm'
$ws.Range("H14").Value = '78%'
$ws.Range("I14").Value = '0.0 mm'
$ws.Range("K14").Value = '0.0 MJ/m2'
$ws.Range("L14").Value = '55.8 km/h - 202º 0:18 TU'
$ws.Range("M14").Value = '-3.2 °C 0:00 TU'
$ws.Range("N14").Value = '-3.3 °C 0:17 TU'
$ws.Range("O14").Value = '-3.2 °C'

# Row 15
$ws.Range("E15").Value = '2026-02-06 00:48:36'
$ws.Range("H15").Value = '56%'
$ws.Range("I15").Value = '0.0 mm'
$ws.Range("J15").Value = '991.2 hPa'
$ws.Range("K15").Value = '0.0 MJ/m2'
$ws.Range("L15").Value = '14.8 km/h - 256º 0:07 TU'
$ws.Range("M15").Value = '13.8 °C 0:02 TU'
$ws.Range("N15").Value = '12.2 °C 0:29 TU'
$ws.Range("O15").Value = '13.2 °C'

# Row 16
$ws.Range("E16").Value = '2026-02-06 00:48:39'
$ws.Range("H16").Value = '95%'
$ws.Range("I16").Value = '0.0 mm'
$ws.Range("K16").Value = '0.0 MJ/m2'
$ws.Range("L16").Value = '13.0 km/h - 306º 0:19 TU'
$ws.Range("M16").Value = '4.4 °C 0:26 TU'
$ws.Range("N16").Value = '4.1 °C 0:15 TU'
$ws.Range("O16").Value = '4.2 °C'

# Row 17
$ws.Range("E17").Value = '2026-02-06 00:48:41'
$ws.Range("H17").Value = '100%'
$ws.Range("I17").Value = '0.0 mm'
$ws.Range("J17").Value = '995.9 hPa'
$ws.Range("K17").Value = '0.0 MJ/m2'
$ws.Range("L17").Value = '10.8 km/h - 215º 0:21 TU'
$ws.Range("M17").Value = '3.6 °C 0:19 TU'
$ws.Range("N17").Value = '3.1 °C 0:24 TU'
$ws.Range("O17").Value = '3.3 °C'

# Row 18
$ws.Range("E18").Value = '2026-02-06 00:48:44'
$ws.Range("G18").Value = '118 cm'
$ws.Range("H18").Value = '96%'
$ws.Range("I18").Value = '0.0 mm'
$ws.Range("K18").Value = '0.0 MJ/m2'
$ws.Range("L18").Value = '22.7 km/h - 313º 0:25 TU'
$ws.Range("M18").Value = '-4.1 °C 0:05 TU'
$ws.Range("N18").Value = '-4.4 °C 0:27 TU'
$ws.Range("O18").Value = '-4.2 °C'

# Row 19
$ws.Range("E19").Value = '2026-02-06 00:48:47'
$ws.Range("H19").Value = '100%'
$ws.Range("I19").Value = '0.0 mm'
$ws.Range("J19").Value = '996.1 hPa'
$ws.Range("K19").Value = '0.0 MJ/m2'
$ws.Range("L19").Value = '11.2 km/h - 301º 0:27 TU'
$ws.Range("M19").Value = '7.2 °C 0:29 TU'
$ws.Range("N19").Value = '6.8 °C 0:16 TU'
$ws.Range("O19").Value = '6.9 °C'

# Row 20
$ws.Range("E20").Value = '2026-02-06 00:48:49'
$ws.Range("G20").Value = '119 cm'
$ws.Range("H20").Value = '75%'
$ws.Range("I20").Value = '0.0 mm'
$ws.Range("K20").Value = '0.0 MJ/m2'
$ws.Range("L20").Value = '20.5 km/h - 286º 0:28 TU'
$ws.Range("M20").Value = '-1.6 °C 0:29 TU'
$ws.Range("N20").Value = '-2.1 °C 0:06 TU'
$ws.Range("O20").Value = '-1.8 °C'

# Row 21
$ws.Range("E21").Value = '2026-02-06 00:48:52'

# Row 22
$ws.Range("E22").Value = '2026-02-06 00:48:54'
$ws.Range("H22").Value = '62%'
$ws.Range("I22").Value = '0.0 mm'
$ws.Range("K22").Value = '0.0 MJ/m2'
$ws.Range("L22").Value = '20.2 km/h - 265º 0:22 TU'
$ws.Range("M22").Value = '13.7 °C 0:09 TU'
$ws.Range("N22").Value = '11.7 °C 0:00 TU'
$ws.Range("O22").Value = '12.8 °C'

# Row 23
$ws.Range("E23").Value = '2026-02-06 00:48:57'
$ws.Range("H23").Value = '91%'
$ws.Range("I23").Value = '0.0 mm'
$ws.Range("J23").Value = '991.9 hPa'
$ws.Range("K23").Value = '0.0 MJ/m2'
$ws.Range("L23").Value = '13.0 km/h - 47º 0:25 TU'
$ws.Range("M23").Value = '7.6 °C 0:25 TU'
$ws.Range("N23").Value = '7.1 °C 0:00 TU'
$ws.Range("O23").Value = '7.4 °C'

# Row 24
$ws.Range("E24").Value = '2026-02-06 00:49:00'
$ws.Range("H24").Value = '62%'
$ws.Range("I24").Value = '0.0 mm'
$ws.Range("J24").Value = '991.1 hPa'
$ws.Range("K24").Value = '0.0 MJ/m2'
$ws.Range("L24").Value = '7.9 km/h - 27º 0:03 TU'
$ws.Range("M24").Value = '12.7 °C 0:11 TU'
$ws.Range("N24").Value = '12.2 °C 0:00 TU'
$ws.Range("O24").Value = '12.5 °C'

# Row 25
$ws.Range("E25").Value = '2026-02-06 00:49:02'
$ws.Range("H25").Value = '90%'
$ws.Range("I25").Value = '0.0 mm'
$ws.Range("J25").Value = '994.5 hPa'
$ws.Range("K25").Value = '0.0 MJ/m2'
$ws.Range("L25").Value = '10.1 km/h - 288º 0:03 TU'
$ws.Range("M25").Value = '3.1 °C 0:06 TU'
$ws.Range("N25").Value = '2.1 °C 0:29 TU'
$ws.Range("O25").Value = '2.6 °C'

# Row 26
$ws.Range("E26").Value = '2026-02-06 00:49:05'
$ws.Range("G26").Value = '112 cm'
$ws.Range("H26").Value = '82%'
$ws.Range("I26").Value = '0.0 mm'
$ws.Range("K26").Value = '0.0 MJ/m2'
$ws.Range("L26").Value = '20.2 km/h - 27º 0:01 TU'
$ws.Range("M26").Value = '0.3 °C 0:03 TU'
$ws.Range("N26").Value = '0.1 °C 0:17 TU'
$ws.Range("O26").Value = '0.1 °C'

# Row 27
$ws.Range("E27").Value = '2026-02-06 00:49:08'
$ws.Range("H27").Value = '86%'
$ws.Range("I27").Value = '0.0 mm'
$ws.Range("J27").Value = '991.4 hPa'
$ws.Range("K27").Value = '0.0 MJ/m2'
$ws.Range("L27").Value = '12.2 km/h - 73º 0:03 TU'
$ws.Range("M27").Value = '10.7 °C 0:06 TU'
$ws.Range("N27").Value = '9.0 °C 0:14 TU'
$ws.Range("O27").Value = '9.9 °C'

# Row 28
$ws.Range("E28").Value = '2026-02-06 00:49:10'
$ws.Range("H28").Value = '79%'
$ws.Range("I28").Value = '0.0 mm'
$ws.Range("J28").Value = '993.1 hPa'
$ws.Range("L28").Value = '31.3 km/h - 242º 0:05 TU'
$ws.Range("M28").Value = '6.9 °C 0:13 TU'
$ws.Range("N28").Value = '4.9 °C 0:03 TU'
$ws.Range("O28").Value = '5.9 °C'

# Row 29
$ws.Range("E29").Value = '2026-02-06 00:49:13'
$ws.Range("H29").Value = '51%'
$ws.Range("I29").Value = '0.0 mm'
$ws.Range("K29").Value = '0.0 MJ/m2'
$ws.Range("L29").Value = '43.9 km/h - 237º 0:21 TU'
$ws.Range("M29").Value = '14.7 °C 0:19 TU'
$ws.Range("N29").Value = '14.0 °C 0:13 TU'
$ws.Range("O29").Value = '14.4 °C'

# Row 30
$ws.Range("E30").Value = '2026-02-06 00:49:16'
$ws.Range("G30").Value = '57 cm'
$ws.Range("H30").Value = '66%'
$ws.Range("I30").Value = '0.0 mm'
$ws.Range("K30").Value = '0.0 MJ/m2'
$ws.Range("L30").Value = '28.1 km/h - 324º 0:02 TU'
$ws.Range("M30").Value = '-1.6 °C 0:00 TU'
$ws.Range("N30").Value = '-2.0 °C 0:23 TU'
$ws.Range("O30").Value = '-1.8 °C'

# Row 31
$ws.Range("E31").Value = '2026-02-06 00:49:18'
$ws.Range("G31").Value = '1 cm'
$ws.Range("H31").Value = '100%'
$ws.Range("I31").Value = '0.0 mm'
$ws.Range("J31").Value = '995.8 hPa'
$ws.Range("M31").Value = '5.2 °C 0:27 TU'
$ws.Range("N31").Value = '4.9 °C 0:15 TU'
$ws.Range("O31").Value = '5.0 °C'

# Row 32
$ws.Range("E32").Value = '2026-02-06 00:49:21'
$ws.Range("H32").Value = '49%'
$ws.Range("I32").Value = '0.0 mm'
$ws.Range("J32").Value = '993.8 hPa'
$ws.Range("K32").Value = '0.0 MJ/m2'
$ws.Range("L32").Value = '49.3 km/h - 278º 0:06 TU'
$ws.Range("M32").Value = '15.9 °C 0:08 TU'
$ws.Range("N32").Value = '15.4 °C 0:22 TU'
$ws.Range("O32").Value = '15.7 °C'

# Row 33
$ws.Range("E33").Value = '2026-02-06 00:49:24'
$ws.Range("H33").Value = '90%'
$ws.Range("I33").Value = '0.0 mm'
$ws.Range("M33").Value = '9.2 °C 0:01 TU'
$ws.Range("N33").Value = '7.5 °C 0:29 TU'
$ws.Range("O33").Value = '8.4 °C'

# Row 34
$ws.Range("E34").Value = '2026-02-06 00:49:26'
$ws.Range("H34").Value = '72%'
$ws.Range("I34").Value = '0.0 mm'
$ws.Range("K34").Value = '0.0 MJ/m2'
$ws.Range("L34").Value = '34.2 km/h - 323º 0:00 TU'
$ws.Range("M34").Value = '9.9 °C 0:24 TU'
$ws.Range("N34").Value = '9.6 °C 0:07 TU'
$ws.Range("O34").Value = '9.7 °C'

# Row 35
$ws.Range("E35").Value = '2026-02-06 00:49:28'
$ws.Range("G35").Value = '197 cm'
$ws.Range("H35").Value = '96%'
$ws.Range("I35").Value = '0.0 mm'
$ws.Range("K35").Value = '0.0 MJ/m2'
$ws.Range("L35").Value = '0.0 km/h - 0º 0:00 TU'
$ws.Range("M35").Value = '-2.6 °C 0:06 TU'
$ws.Range("N35").Value = '-2.7 °C 0:24 TU'
$ws.Range("O35").Value = '-2.7 °C'

# Row 36
$ws.Range("E36").Value = '2026-02-06 00:49:31'
$ws.Range("H36").Value = '63%'
$ws.Range("I36").Value = '0.0 mm'
$ws.Range("J36").Value = '994.8 hPa'
$ws.Range("K36").Value = '0.0 MJ/m2'
$ws.Range("L36").Value = '30.2 km/h - 2º 0:29 TU'
$ws.Range("M36").Value = '13.2 °C 0:11 TU'
$ws.Range("N36").Value = '12.6 °C 0:04 TU'
$ws.Range("O36").Value = '12.9 °C'
